# Apply updated crypto price/volume data (and two coin-row reorderings)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.964.90"
$ws.Range("E2").Value = "  +1.84%  "
$ws.Range("D3").Value = "1.703.16"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'315.78"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.3982"
$ws.Range("E7").Value = "  +1.54%  "
$ws.Range("D8").Value = "'0.4037"
$ws.Range("E8").Value = "  -0.81%  "
$ws.Range("D9").Value = "'1.469"
$ws.Range("E9").Value = "  -1.54%  "
$ws.Range("D10").Value = "'53.44"
$ws.Range("E10").Value = "  +1.82%  "
$ws.Range("D11").Value = "'1.001"
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").Value = "'0.08809"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").Value = "'25.97"
$ws.Range("E13").Value = "  -2.52%  "
$ws.Range("D14").Value = "'7.464"
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").Value = "'7.977"
$ws.Range("E15").Value = "  -1.98%  "
$ws.Range("D16").Value = "'0.00001351"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").Value = "1.774.22"
$ws.Range("E17").Value = "  +5.19%  "
$ws.Range("D18").Value = "'96.04"
$ws.Range("E18").Value = "  -2.08%  "
$ws.Range("D19").Value = "'0.07198"
$ws.Range("E19").Value = "  +0.69%  "
$ws.Range("D20").Value = "'20.67"
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("D21").Value = "'7.308"
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").Value = "'0.9996"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").Value = "'14.32"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").Value = "24.963.97"
$ws.Range("E24").Value = "  +1.88%  "
$ws.Range("D25").Value = "'2.402"
$ws.Range("E25").Value = "  +3.41%  "
$ws.Range("D26").Value = "'2.946"
$ws.Range("E26").Value = "  -2.09%  "
$ws.Range("D27").Value = "'23.56"
$ws.Range("E27").Value = "  +3.66%  "
$ws.Range("D28").Value = "'6.113"
$ws.Range("E28").Value = "  +13.26%  "
$ws.Range("D29").Value = "'162.57"
$ws.Range("E29").Value = "  -2.51%  "
$ws.Range("D30").Value = "'151.03"
$ws.Range("E30").Value = "  +8.35%  "
$ws.Range("D31").Value = "'8.462"
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("D32").Value = "'2.670"
$ws.Range("E32").Value = "  +21.62%  "
$ws.Range("D33").Value = "1.952.54"
$ws.Range("E33").Value = "  +4.33%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.08548"
$ws.Range("E34").Value = "  -2.60%  "
$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").Value = "'0.03161"
$ws.Range("E35").Value = "  +5.70%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.046"
$ws.Range("E36").Value = "  +0.77%  "
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "'7.208"
$ws.Range("E37").Value = "  -1.68%  "
$ws.Range("D38").Value = "'0.2883"
$ws.Range("E38").Value = "  +3.55%  "
$ws.Range("D39").Value = "'10.95"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").Value = "'0.09585"
$ws.Range("E40").Value = "  +4.64%  "
$ws.Range("D41").Value = "'0.8248"
$ws.Range("E41").Value = "  +2.57%  "
$ws.Range("D42").Value = "'14.01"
$ws.Range("E42").Value = "  -1.33%  "
$ws.Range("D43").Value = "'1.480"
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("D44").Value = "'17.18"
$ws.Range("E44").Value = "  -2.27%  "
$ws.Range("D45").Value = "'2.684"
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("D46").Value = "'0.7389"
$ws.Range("E46").Value = "  +1.55%  "
$ws.Range("D47").Value = "'4.253"
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.08887"
$ws.Range("E48").Value = "  +8.85%  "
$ws.Range("B49").Value = "Flow"
$ws.Range("C49").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D49").Value = "'1.399"
$ws.Range("E49").Value = "  -0.46%  "
$ws.Range("D50").Value = "'0.9999"
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("D51").Value = "'139.52"
$ws.Range("E51").Value = "  -1.10%  "
